$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of VIN data (row 6) - Volkswagen Arteon SEL
$ws.Range("B6").Value = "SYMBOL_2000_CHOICE_T"
$ws.Range("C6").Value = 2017
$ws.Range("D6").Value = "Volkswagen"
$ws.Range("E6").Value = "Volkswagen"
$ws.Range("F6").Value = "Arteon"
$ws.Range("G6").Value = "Arteon SEL"
$ws.Range("H6").Value = 88888
$ws.Range("I6").Value = "WAG"
$ws.Range("J6").Value = "Coupe"
$ws.Range("K6").Value = "Sedan"
$ws.Range("L6").Value = "Coupe"
$ws.Range("M6").Value = "WAG"
$ws.Range("N6").Value = "8L V12"
$ws.Range("O6").Value = 12
$ws.Range("P6").Value = "G"
$ws.Range("Q6").Value = 214
$ws.Range("R6").Value = "4WD"
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = "000R"
$ws.Range("U6").Value = "DUAL AIR BAGS FRONT"
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = "4 WHEEL STANDARD"
$ws.Range("X6").Value = "STD"
$ws.Range("Y6").Value = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
$ws.Range("Z6").Value = "GC"
$ws.Range("AA6").Value = 35
$ws.Range("AB6").Value = 44
$ws.Range("AC6").Value = "S"
$ws.Range("AD6").Value = "Y"
$ws.Range("AE6").Value = "R"
$ws.Range("AF6").Value = "E"
$ws.Range("AG6").Value = "S"
$ws.Range("AH6").Value = "A"
$ws.Range("AI6").Value = 20000101
$ws.Range("AJ6").Value = "Y"
$ws.Range("AK6").Value = "Y"
$ws.Range("AL6").Value = "N"
$ws.Range("A6").Value = "8MSRP17H&V"

# Match original formatting: B6 keeps its pre-existing style (the "Good" green style),
# the rest of the row (C6:AL6) use the left-aligned style used by the other data rows.
$ws.Range("C6:AL6").Style = "Normal"
$ws.Range("C6:AL6").HorizontalAlignment = -4131

# Update selection to match the saved view state
$ws.Range("A8").Select()
